$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 790.6667
$ws.Range("J9").Value = 343.5
$ws.Range("L9").Value = 343.5
$ws.Range("N9").Value = -681.5
$ws.Range("H80").Value = 8964919
$ws.Range("J80").Value = 50987.875
$ws.Range("L80").Value = 152963.625
$ws.Range("N80").Value = -154959.625
$ws.Range("H83").Value = 8964919
$ws.Range("J83").Value = 50987.875
$ws.Range("L83").Value = 458890.875
$ws.Range("N83").Value = -468874.875
$ws.Range("H132").Value = 979.3103599999999
$ws.Range("I132").Value = 856.3200000000001
$ws.Range("J132").Value = 1748
$ws.Range("K132").Value = 2568.96
$ws.Range("L132").Value = 5244
$ws.Range("M132").Value = -38.96000000000004
$ws.Range("N132").Value = -10304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H42").Value = 24970
$ws.Range("J42").Value = 24970
$ws.Range("L42").Value = 24970
$ws.Range("N42").Value = -25942
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H63").Value = 2223.25
$ws.Range("I63").Value = 1966
$ws.Range("K63").Value = 1966
$ws.Range("M63").Value = -1280
$ws.Range("H66").Value = 2223.25
$ws.Range("I66").Value = 1966
$ws.Range("K66").Value = 9830
$ws.Range("M66").Value = -6398
$ws.Range("H74").Value = 23906.088
$ws.Range("I74").Value = 27727.71
$ws.Range("J74").Value = 5753.375
$ws.Range("K74").Value = 27727.71
$ws.Range("L74").Value = 5753.375
$ws.Range("M74").Value = -26853.71
$ws.Range("N74").Value = -7501.375
$ws.Range("H77").Value = 23906.088
$ws.Range("I77").Value = 27727.71
$ws.Range("J77").Value = 5753.375
$ws.Range("K77").Value = 138638.55
$ws.Range("L77").Value = 28766.875
$ws.Range("M77").Value = -134270.55
$ws.Range("N77").Value = -37502.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5750134.5
$ws.Range("I20").Value = 9806538
$ws.Range("J20").Value = 3562.1667
$ws.Range("K20").Value = 9806538
$ws.Range("L20").Value = 3562.1667
$ws.Range("M20").Value = -9806291
$ws.Range("N20").Value = -4056.1667
$ws.Range("H94").Value = 1892.6111
$ws.Range("I94").Value = 760.7692
$ws.Range("K94").Value = 760.7692
$ws.Range("M94").Value = -309.7692
$ws.Range("H103").Value = 17752.166
$ws.Range("J103").Value = 17752.166
$ws.Range("L103").Value = 17752.166
$ws.Range("N103").Value = -20096.166
$ws.Range("H134").Value = 5211905
$ws.Range("I134").Value = 8622133
$ws.Range("J134").Value = 6820.684
$ws.Range("K134").Value = 25866399
$ws.Range("L134").Value = 20462.052
$ws.Range("M134").Value = -25863864
$ws.Range("N134").Value = -25532.052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 841
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -450
$ws.Range("H59").Value = 50664
$ws.Range("J59").Value = 50664
$ws.Range("L59").Value = 50664
$ws.Range("N59").Value = -52954
$ws.Range("H94").Value = 1296.6666
$ws.Range("I94").Value = 1324.625
$ws.Range("K94").Value = 1324.625
$ws.Range("M94").Value = -873.625
$ws.Range("H105").Value = 11906304
$ws.Range("I105").Value = 14286565
$ws.Range("K105").Value = 14286565
$ws.Range("M105").Value = -14284818
$ws.Range("H132").Value = 2757.309
$ws.Range("I132").Value = 1929.975
$ws.Range("K132").Value = 5789.924999999999
$ws.Range("M132").Value = -3259.924999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 30982.25
$ws.Range("I120").Value = 7930
$ws.Range("J120").Value = 38666.332
$ws.Range("K120").Value = 23790
$ws.Range("L120").Value = 115998.996
$ws.Range("M120").Value = -18952
$ws.Range("N120").Value = -125674.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H52").Value = 77005.8
$ws.Range("J52").Value = 90000
$ws.Range("L52").Value = 90000
$ws.Range("N52").Value = -90518
$ws.Range("H58").Value = 66323.5
$ws.Range("J58").Value = 78980
$ws.Range("L58").Value = 78980
$ws.Range("N58").Value = -79534
$ws.Range("H64").Value = 74656.5
$ws.Range("H67").Value = 74656.5
$ws.Range("H70").Value = 7171.2856
$ws.Range("I70").Value = 4762.857
$ws.Range("J70").Value = 8375.5
$ws.Range("K70").Value = 4762.857
$ws.Range("L70").Value = 8375.5
$ws.Range("M70").Value = -4492.857
$ws.Range("N70").Value = -8915.5
$ws.Range("H73").Value = 7171.2856
$ws.Range("I73").Value = 4762.857
$ws.Range("J73").Value = 8375.5
$ws.Range("K73").Value = 4762.857
$ws.Range("L73").Value = 8375.5
$ws.Range("M73").Value = -3826.857
$ws.Range("N73").Value = -10247.5
$ws.Range("H102").Value = 3245.45
$ws.Range("I102").Value = 3088.162
$ws.Range("K102").Value = 3088.162
$ws.Range("M102").Value = -1466.162
$ws.Range("H139").Value = 66665.60000000001
$ws.Range("J139").Value = 66665.60000000001
$ws.Range("L139").Value = 66665.60000000001
$ws.Range("N139").Value = -76945.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 16678331
$ws.Range("H69").Value = 53863
$ws.Range("J69").Value = 53863
$ws.Range("L69").Value = 53863
$ws.Range("N69").Value = -55485
$ws.Range("H72").Value = 53863
$ws.Range("J72").Value = 53863
$ws.Range("L72").Value = 161589
$ws.Range("N72").Value = -169701
$ws.Range("H101").Value = 69750
$ws.Range("J101").Value = 69750
$ws.Range("L101").Value = 69750
$ws.Range("N101").Value = -76240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 30072570
$ws.Range("I81").Value = 1751332.1
$ws.Range("K81").Value = 3502664.2
$ws.Range("M81").Value = -3501603.2
$ws.Range("H84").Value = 30072570
$ws.Range("I84").Value = 1751332.1
$ws.Range("K84").Value = 17513321
$ws.Range("M84").Value = -17508017
$ws.Range("H132").Value = 1624.7273
$ws.Range("J132").Value = 5277.5
$ws.Range("L132").Value = 15832.5
$ws.Range("N132").Value = -20892.5
